$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "One rainy afternoon, Emily's beloved teddy bear went missing. Frantically, she searched every corner of her room, but the bear was nowhere to be found. With tears in her eyes, Emily ventured outside into the downpour, determined to find her lost friend. As she wandered through the streets, she spotted a familiar sight - a soggy teddy bear lying abandoned on a bench. With a joyful cry, Emily rushed forward and scooped up her precious bear. Grateful and relieved, she hugged her teddy tightly, knowing that some bonds are unbreakable, even in the stormiest of times."

$ws.Range("C2").Value = "Em uma tarde chuvosa, o amado urso de pelúcia de Emily desapareceu. Desesperadamente, ela procurou em todos os cantos do seu quarto, mas o urso não estava em lugar nenhum. Com lágrimas nos olhos, Emily aventurou-se do lado de fora sob a chuva forte, determinada a encontrar seu amigo perdido. Enquanto vagava pelas ruas, ela avistou uma visão familiar - um urso de pelúcia encharcado deitado abandonado em um banco. Com um grito alegre, Emily correu para frente e pegou seu urso precioso. Grata e aliviada, ela abraçou seu urso com força, sabendo que alguns laços são inquebráveis, mesmo nos tempos mais tempestuosos."

$ws.Range("D2").Value = "Básico Avançado"

$ws.Range("E2").Value = 250
